$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy number formats (date/time styles) from the last existing row (44) down to the
# new rows (45-66) for columns A:M so the new cells pick up the same cell styles
# (s="1" for dates in col A, s="2" for times in col B) as the rest of the sheet.
$ws.Range("A44:M44").Copy()
$ws.Range("A45:M66").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 45 col B uses the h:mm (no AM/PM) style like row 4/5/35/43, so pull that format in specifically.
$ws.Range("B43").Copy()
$ws.Range("B45").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A45").Value = 42104
$ws.Range("B45").Value = 0.5833333333333334
$ws.Range("C45").Value = "RP"
$ws.Range("D45").Value = 425
$ws.Range("E45").Value = 4
$ws.Range("F45").Value = 1
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 446
$ws.Range("I45").Value = 22
$ws.Range("J45").Value = 2146
$ws.Range("K45").Value = 87
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = "maybe"
$ws.Range("N45").Value = "said near the end he was a little suspicious of the essay evaluations"

$ws.Range("A46").Value = 42104
$ws.Range("B46").Value = 0.5833333333333334
$ws.Range("C46").Value = "VP "
$ws.Range("D46").Value = 426
$ws.Range("E46").Value = 1
$ws.Range("F46").Value = 2
$ws.Range("G46").Value = 1
$ws.Range("H46").Value = 202
$ws.Range("I46").Value = 4
$ws.Range("J46").Value = 797
$ws.Range("K46").Value = 85
$ws.Range("L46").Value = 56
$ws.Range("M46").Value = "yes "
$ws.Range("N46").Value = "subject was quiet "

$ws.Range("A47").Value = 42104
$ws.Range("B47").Value = 0.625
$ws.Range("C47").Value = "RP"
$ws.Range("D47").Value = 427
$ws.Range("E47").Value = 4
$ws.Range("F47").Value = 3
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 603
$ws.Range("I47").Value = 20
$ws.Range("J47").Value = 49
$ws.Range("K47").Value = 515
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = "yes "

$ws.Range("A48").Value = 42104
$ws.Range("B48").Value = 0.625
$ws.Range("C48").Value = "VP "
$ws.Range("D48").Value = 428
$ws.Range("E48").Value = 1
$ws.Range("F48").Value = 4
$ws.Range("G48").Value = 14
$ws.Range("H48").Value = 232
$ws.Range("I48").Value = 3
$ws.Range("J48").Value = 214
$ws.Range("K48").Value = 138
$ws.Range("L48").Value = 181
$ws.Range("M48").Value = "yes "
$ws.Range("N48").Value = "participant died on 1st level twice "

$ws.Range("A49").Value = 42104
$ws.Range("B49").Value = 0.6666666666666666
$ws.Range("C49").Value = "RP"
$ws.Range("D49").Value = 429
$ws.Range("E49").Value = 4
$ws.Range("F49").Value = 1
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 206
$ws.Range("I49").Value = 7
$ws.Range("J49").Value = 37
$ws.Range("K49").Value = 219
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = "maybe"
$ws.Range("N49").Value = "said the evaluation seemed to harsh to be real "

$ws.Range("A50").Value = 42109
$ws.Range("B50").Value = 0.5833333333333334
$ws.Range("C50").Value = "TM"
$ws.Range("D50").Value = 430
$ws.Range("E50").Value = 4
$ws.Range("F50").Value = 2
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 311
$ws.Range("I50").Value = 10
$ws.Range("J50").Value = 1118
$ws.Range("K50").Value = 70
$ws.Range("L50").Value = 47
$ws.Range("M50").Value = "maybe"
$ws.Range("N50").Value = "felt like ther may have been a fake partner after reading insults "

$ws.Range("A51").Value = 42109
$ws.Range("B51").Value = 0.625
$ws.Range("C51").Value = "TM"
$ws.Range("D51").Value = 431
$ws.Range("E51").Value = 4
$ws.Range("F51").Value = 3
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 471
$ws.Range("I51").Value = 16
$ws.Range("J51").Value = 515
$ws.Range("K51").Value = 322
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = "yes"

$ws.Range("A52").Value = 42109
$ws.Range("B52").Value = 0.625
$ws.Range("C52").Value = "AH "
$ws.Range("D52").Value = 432
$ws.Range("E52").Value = 1
$ws.Range("F52").Value = 4
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 411
$ws.Range("I52").Value = 13
$ws.Range("J52").Value = 1450
$ws.Range("K52").Value = 127
$ws.Range("L52").Value = 84
$ws.Range("M52").Value = "yes"
$ws.Range("N52").Value = "showed dissapointment when told we werent going to continue with the distraction task "

$ws.Range("A53").Value = 42109
$ws.Range("B53").Value = 0.6666666666666666
$ws.Range("C53").Value = "AH "
$ws.Range("D53").Value = 433
$ws.Range("E53").Value = 1
$ws.Range("F53").Value = 1
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 436
$ws.Range("I53").Value = 14
$ws.Range("J53").Value = 2614
$ws.Range("K53").Value = 6
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = "yes"

$ws.Range("A54").Value = 42109
$ws.Range("B54").Value = 0.6666666666666666
$ws.Range("C54").Value = "TM"
$ws.Range("D54").Value = 434
$ws.Range("E54").Value = 4
$ws.Range("F54").Value = 2
$ws.Range("G54").Value = 4
$ws.Range("H54").Value = 383
$ws.Range("I54").Value = 6
$ws.Range("J54").Value = 1464
$ws.Range("K54").Value = 92
$ws.Range("L54").Value = 126
$ws.Range("M54").Value = "yes "

$ws.Range("A55").Value = 42111
$ws.Range("B55").Value = 0.625
$ws.Range("C55").Value = "RP"
$ws.Range("D55").Value = 435
$ws.Range("E55").Value = 2
$ws.Range("F55").Value = 3
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 459
$ws.Range("I55").Value = 15
$ws.Range("J55").Value = 2240
$ws.Range("K55").Value = 56
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = "no"
$ws.Range("N55").Value = "said after the evaluation he knew that it was fake and if he had picked pro life he would've got a pro-choice essay back because they were probably fake "

$ws.Range("A56").Value = 42111
$ws.Range("B56").Value = 0.625
$ws.Range("C56").Value = "VP "
$ws.Range("D56").Value = 436
$ws.Range("E56").Value = 4
$ws.Range("F56").Value = 4
$ws.Range("G56").Value = 6
$ws.Range("H56").Value = 288
$ws.Range("I56").Value = 3
$ws.Range("J56").Value = 1181
$ws.Range("K56").Value = 72
$ws.Range("L56").Value = 116
$ws.Range("M56").Value = "maybe"
$ws.Range("N56").Value = "subject stated ge was suspicious tat the iinsult evaluation was a fake during funneled debriefing. He said it was because the evaluation was unrealistically harsh "

$ws.Range("A57").Value = 42111
$ws.Range("B57").Value = 0.6666666666666666
$ws.Range("C57").Value = "RP"
$ws.Range("D57").Value = 437
$ws.Range("E57").Value = 4
$ws.Range("F57").Value = 1
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 485
$ws.Range("I57").Value = 22
$ws.Range("J57").Value = 670
$ws.Range("K57").Value = 240
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = "maybe"
$ws.Range("N57").Value = "he didn’t explicitly say that he knew what was going on but he didn’t seem surprised at all and did say the evaluation seemed fishy "

$ws.Range("A58").Value = 42116
$ws.Range("B58").Value = 0.625
$ws.Range("C58").Value = "TM"
$ws.Range("D58").Value = 438
$ws.Range("E58").Value = 4
$ws.Range("F58").Value = 2
$ws.Range("G58").Value = 2
$ws.Range("H58").Value = 263
$ws.Range("I58").Value = 3
$ws.Range("J58").Value = 1031
$ws.Range("K58").Value = 52
$ws.Range("L58").Value = 64
$ws.Range("M58").Value = "maybe"
$ws.Range("N58").Value = "said it was on of the most violent games hes played; said insults were super mean, so he got a little skeptical "

$ws.Range("A59").Value = 42116
$ws.Range("B59").Value = 0.6666666666666666
$ws.Range("C59").Value = "TM"
$ws.Range("D59").Value = 439
$ws.Range("E59").Value = 4
$ws.Range("F59").Value = 3
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 425
$ws.Range("I59").Value = 14
$ws.Range("J59").Value = 681
$ws.Range("K59").Value = 310
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = "maybe "
$ws.Range("N59").Value = "said insults were a;; negatives so either the other guy was really mean or we were trying to pull something "

$ws.Range("A60").Value = 42116
$ws.Range("B60").Value = 0.6666666666666666
$ws.Range("C60").Value = "AH "
$ws.Range("D60").Value = 440
$ws.Range("E60").Value = 1
$ws.Range("F60").Value = 4
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 385
$ws.Range("I60").Value = 13
$ws.Range("J60").Value = 1958
$ws.Range("K60").Value = 77
$ws.Range("L60").Value = 78
$ws.Range("M60").Value = "maybe"
$ws.Range("N60").Value = "said he had a few suspicions after skipping the distraction assignment even after testing the water "

$ws.Range("A61").Value = 42118
$ws.Range("B61").Value = 0.625
$ws.Range("C61").Value = "RP"
$ws.Range("D61").Value = 441
$ws.Range("E61").Value = 1
$ws.Range("F61").Value = 1
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 444
$ws.Range("I61").Value = 22
$ws.Range("J61").Value = 612
$ws.Range("K61").Value = 264
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = "yes "
$ws.Range("N61").Value = "said that the evaluation seemed kind of weird. After being debriefed the subject said he was genuinely surprised about everything and `"we had got him`" "

$ws.Range("A62").Value = 42123
$ws.Range("B62").Value = 0.625
$ws.Range("C62").Value = "TM"
$ws.Range("D62").Value = 442
$ws.Range("E62").Value = 4
$ws.Range("F62").Value = 2
$ws.Range("G62").Value = 1
$ws.Range("H62").Value = 373
$ws.Range("I62").Value = 12
$ws.Range("J62").Value = 1419
$ws.Range("K62").Value = 124
$ws.Range("L62").Value = 56
$ws.Range("M62").Value = "yes"

$ws.Range("A63").Value = 42128
$ws.Range("B63").Value = 0.625
$ws.Range("C63").Value = "TM "
$ws.Range("D63").Value = 443
$ws.Range("E63").Value = 4
$ws.Range("F63").Value = 3
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 392
$ws.Range("I63").Value = 17
$ws.Range("J63").Value = 991
$ws.Range("K63").Value = 183
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = "yes"

$ws.Range("A64").Value = 42128
$ws.Range("B64").Value = 0.6666666666666666
$ws.Range("C64").Value = "TM"
$ws.Range("D64").Value = 444
$ws.Range("E64").Value = 4
$ws.Range("F64").Value = 4
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = 440
$ws.Range("I64").Value = 12
$ws.Range("J64").Value = 1048
$ws.Range("K64").Value = 176
$ws.Range("L64").Value = 89
$ws.Range("M64").Value = "maybe "
$ws.Range("N64").Value = "said insults were really mean and kind of made him think it was fake "

$ws.Range("A65").Value = 42130
$ws.Range("B65").Value = 0.5833333333333334
$ws.Range("C65").Value = "TM"
$ws.Range("D65").Value = 445
$ws.Range("E65").Value = 4
$ws.Range("F65").Value = 1
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 439
$ws.Range("I65").Value = 15
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 368
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = "yes"

$ws.Range("A66").Value = 42130
$ws.Range("B66").Value = 0.625
$ws.Range("C66").Value = "TM"
$ws.Range("D66").Value = 446
$ws.Range("E66").Value = 4
$ws.Range("F66").Value = 2
$ws.Range("G66").Value = 1
$ws.Range("H66").Value = 294
$ws.Range("I66").Value = 9
$ws.Range("J66").Value = 1176
$ws.Range("K66").Value = 72
$ws.Range("L66").Value = 59
$ws.Range("M66").Value = "no"
$ws.Range("N66").Value = "said he's done a study very similar to this, so he assumed he wasn’t paired with anyone "

# Update the view/selection state to match where data entry left off.
$ws.Range("A66").Select()

Write-Output "done"